$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "238.07"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.64"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.457"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05641"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.495"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.354"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.7933"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.055"
$ws.Range("D9").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03203"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02994"
$ws.Range("D13").Style = "Normal"

$ws.Range("B14").Value = "ProBitToken"

$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1054"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "13ProBitTokenPROB"

$ws.Range("B15").Value = "BitMartToken"

$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09244"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"

$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001660"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "MCDex"

$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.254"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "CoinExToken"

$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04771"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("B19").Value = "One"

$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005745"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "18OneONEWorstin24h"

$ws.Range("B20").Value = "TigerCash"

$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006222"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"

$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.005096"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"

$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001052"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"

$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001502"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "UpBots"

$ws.Range("C24").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0003203"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "23UpBotsUBXT"

$ws.Range("B25").Value = "LEO"

$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.912"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "24LEOLEOBestin24h"

$ws.Range("B26").Value = "BTSEToken"

$ws.Range("C26").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.201"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "25BTSETokenBTSE"

$ws.Range("B27").Value = "BitpandaEcosystemToken"

$ws.Range("C27").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.3337"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "26BitpandaEcosystemTokenBEST"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04117"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006925"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003504"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1039"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008758"
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6758"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03608"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "47BOLOBOLO"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
$ws.Range("D50").Style = "Normal"
